$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells keep their exact text representation (avoid Excel auto-numeric conversion)
foreach ($addr in @("D2","D3","D5","D6","D7","D10","D11","D15","D16","D18","D19","D22","D23","D24","D25","D27","D30","D31","D32","D35","D37","D39","D40","D41","D42","D43","D44","D45","D49","D50","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "48.568.96"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").Value = "2.611.20"
$ws.Range("E3").Value = "  +2.48%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "321.15"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").Value = "108.84"
$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("D7").Value = "0.520"
$ws.Range("E7").Value = "  -1.35%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("E9").Value = "  -2.88%  "

$ws.Range("D10").Value = "39.02"
$ws.Range("E10").Value = "  -2.53%  "

$ws.Range("D11").Value = "19.69"
$ws.Range("E11").Value = "  -3.19%  "

$ws.Range("E12").Value = "  -1.12%  "

$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("D15").Value = "3.022.18"
$ws.Range("E15").Value = "  +2.78%  "

$ws.Range("D16").Value = "2.570.98"
$ws.Range("E16").Value = "  +0.96%  "

$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("D18").Value = "48.568.83"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").Value = "2.92"
$ws.Range("E19").Value = "  -3.76%  "

$ws.Range("E20").Value = "  -3.69%  "

$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("D22").Value = "0.0₃0937"
$ws.Range("E22").Value = "  -0.72%  "

$ws.Range("D23").Value = "268.60"
$ws.Range("E23").Value = "  -5.34%  "

$ws.Range("D24").Value = "68.54"
$ws.Range("E24").Value = "  -4.70%  "

$ws.Range("D25").Value = "2.51"
$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("E26").Value = "  -1.79%  "

$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.21%  "

$ws.Range("E28").Value = "  +1.54%  "

$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "34.61"
$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "0.136"
$ws.Range("E31").Value = "  -5.24%  "

$ws.Range("D32").Value = "49.09"

$ws.Range("E33").Value = "  +1.30%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").Value = "18.99"
$ws.Range("E35").Value = "  -3.17%  "

$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").Value = "4.93"
$ws.Range("E37").Value = "  +5.25%  "

$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("D39").Value = "3.11"
$ws.Range("E39").Value = "  +3.78%  "

$ws.Range("D40").Value = "125.14"
$ws.Range("E40").Value = "  +3.02%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").Value = "  -1.65%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "22.17"
$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").Value = "2.12"
$ws.Range("E43").Value = "  -3.98%  "

$ws.Range("D44").Value = "0.0313"
$ws.Range("E44").Value = "  +0.76%  "

$ws.Range("D45").Value = "2.054.67"
$ws.Range("E45").Value = "  +2.02%  "

$ws.Range("E46").Value = "  -4.26%  "

$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("D49").Value = "8.87"
$ws.Range("E49").Value = "  -1.33%  "

$ws.Range("D50").Value = "58.25"
$ws.Range("E50").Value = "  +2.20%  "

$ws.Range("D51").Value = "5.14"
$ws.Range("E51").Value = "  -3.14%  "
